$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (model "linear" -> "xgb") with new metric values
$ws.Range("B2").Value = "xgb"
$ws.Range("C2").Value = 219.34
$ws.Range("D2").Value = 90055.77
$ws.Range("E2").Value = 0.73
$ws.Range("F2").Value = 124.04
$ws.Range("G2").Value = 34550.16
$ws.Range("H2").Value = 0.88
$ws.Range("I2").Value = 19
$ws.Range("J2").Value = 684.76
$ws.Range("K2").Value = 1

# Delete row 3 entirely (the "baseline-rent" row)
$ws.Rows("3:3").Delete()
